# Generate Report for Handoff
# Updates the localization-status workbook so the generated file id
# (2f7d660e-d518-4bec-9a74-221413f83b3d -> e3a8775b-65f1-4368-91a2-95eef1b2486e)
# and the associated handoff/handback timestamps & xliff file names are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# A2: bare markdown file name
$wsOverview.Range("A2").Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
# B2: path-and-name (hyperlinked) -- also refresh the hyperlink display text
$wsOverview.Range("B2").Value = "e2e\e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-01 21:07:00"

# --- zh-cn sheet ------------------------------------------------------
# A2: bare markdown file name (hyperlinked)
$wsZhCn.Range("A2").Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
# G2: Latest Handoff File
$wsZhCn.Range("G2").Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.f50278a990b0bd43566c460405169ea95900c7d9.zh-cn.xlf"
# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-01 21:06:55"

# --- de-de sheet ------------------------------------------------------
# A2: bare markdown file name (hyperlinked)
$wsDeDe.Range("A2").Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
# G2: Latest Handback File
$wsDeDe.Range("G2").Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.f50278a990b0bd43566c460405169ea95900c7d9.de-de.xlf"
# H2: Latest Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-01 21:07:00"
